# "Turklerde Eksik Kalan Haziran Hakedisleri Eklendi"
# Adds the missing June ("HAZIRAN KALAN", column F) carry-over hours for
# three employees (rows 3, 4 and 12), lets the dependent formulas in
# column I (and the shared K/H formulas) recalculate, puts a medium
# right-hand border on the K column totals, and leaves the selection on
# F14 (the last cell the user touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (SADIK ACAR): add the missing June carry-over hours ----------
$ws.Range("F3").Value = 1883
$ws.Range("I3").Formula = "=(F3+H3)*M3"

# --- Row 4 (IHSAN GOL): add the missing June carry-over hours -----------
$ws.Range("F4").Value = 1142
$ws.Range("I4").Formula = "=(F4+H4)*M4"

# --- Row 12 (BAYRAM GONCE): add the missing June carry-over hours -------
$ws.Range("F12").Value = 597

# --- Put a medium right border on the K (NET HAKEDIS) column cells so it
#     stands out from the (now border-less) L column -------------------
$kRange = $ws.Range("K3:K14")
$kRange.Borders.Item(10).LineStyle = 1
$kRange.Borders.Item(10).Weight = -4138

# --- Restore the active selection to F14, matching the saved workbook --
$ws.Range("F14").Select()
